$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: merge "A" + " " runs into "A " (keep "slide" separate) ---
$sh1 = $s.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$titleMerge = $tr1.Characters(1, 2)
$titleMerge.Text = $titleMerge.Text

# --- Table cell (row1,col2): merge "a" + " " runs into "a " (keep "table" separate) ---
$tblShape = $s.Shapes.Item(3)
$tbl = $tblShape.Table
$cellTr = $tbl.Cell(1, 2).Shape.TextFrame.TextRange
$fullLen = $cellTr.Text.Length
$cellMerge = $cellTr.Characters(1, $fullLen)
$cellMerge.Text = "a "
